$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2947.4167
$ws.Range("I32").Value = 3111.25
$ws.Range("J32").Value = 2619.75
$ws.Range("K32").Value = 3111.25
$ws.Range("L32").Value = 2619.75
$ws.Range("M32").Value = -2785.25
$ws.Range("N32").Value = -3271.75

$ws.Range("H94").Value = 2143.75
$ws.Range("I94").Value = 2143.75
$ws.Range("K94").Value = 2143.75
$ws.Range("M94").Value = -1692.75

$ws.Range("H132").Value = 3207.06
$ws.Range("I132").Value = 1405.5
$ws.Range("J132").Value = 7035.375
$ws.Range("K132").Value = 4216.5
$ws.Range("L132").Value = 21106.125
$ws.Range("M132").Value = -1686.5
$ws.Range("N132").Value = -26166.125

$ws.Range("H137").Value = 31746.424
$ws.Range("I137").Value = 1532.7391
$ws.Range("K137").Value = 4598.2173
$ws.Range("M137").Value = -2048.2173

$ws.Range("H138").Value = 2496.4285
$ws.Range("I138").Value = 1589.125
$ws.Range("K138").Value = 4767.375
$ws.Range("M138").Value = 372.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2316.1562
$ws.Range("I2").Value = 2155.261
$ws.Range("K2").Value = 2155.261
$ws.Range("M2").Value = -2042.261

$ws.Range("H102").Value = 1530.3529
$ws.Range("I102").Value = 1472.862
$ws.Range("J102").Value = 1863.8
$ws.Range("K102").Value = 1472.862
$ws.Range("L102").Value = 1863.8
$ws.Range("M102").Value = 149.1379999999999
$ws.Range("N102").Value = -5107.8

$ws.Range("H116").Value = 2316.1562
$ws.Range("I116").Value = 2155.261
$ws.Range("K116").Value = 2155.261
$ws.Range("M116").Value = 138.739

$ws.Range("H132").Value = 24750.5
$ws.Range("I132").Value = 27487.174
$ws.Range("K132").Value = 82461.522
$ws.Range("M132").Value = -79931.522

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2316.1562
$ws.Range("I3").Value = 2155.261
$ws.Range("K3").Value = 2155.261
$ws.Range("M3").Value = -2041.261

$ws.Range("H82").Value = 22248.555
$ws.Range("J82").Value = 24997.143
$ws.Range("L82").Value = 24997.143
$ws.Range("N82").Value = -25763.143

$ws.Range("H85").Value = 22248.555
$ws.Range("J85").Value = 24997.143
$ws.Range("L85").Value = 24997.143
$ws.Range("N85").Value = -27649.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2119.5264
$ws.Range("I31").Value = 2204.2942
$ws.Range("K31").Value = 2204.2942
$ws.Range("M31").Value = -1909.2942

$ws.Range("H34").Value = 2119.5264
$ws.Range("I34").Value = 2204.2942
$ws.Range("K34").Value = 2204.2942
$ws.Range("M34").Value = -2002.2942

$ws.Range("H50").Value = 11110.556
$ws.Range("J50").Value = 11110.556
$ws.Range("L50").Value = 11110.556
$ws.Range("N50").Value = -12360.556

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H59").Value = 15998.8
$ws.Range("J59").Value = 15998.8
$ws.Range("L59").Value = 15998.8
$ws.Range("N59").Value = -18288.8

$ws.Range("H60").Value = 10672.333
$ws.Range("J60").Value = 10768.692
$ws.Range("L60").Value = 10768.692
$ws.Range("N60").Value = -11790.692

$ws.Range("H68").Value = 24996.666
$ws.Range("J68").Value = 24996.666
$ws.Range("L68").Value = 24996.666
$ws.Range("N68").Value = -26494.666

$ws.Range("H71").Value = 24996.666
$ws.Range("J71").Value = 24996.666
$ws.Range("L71").Value = 74989.99800000001
$ws.Range("N71").Value = -82477.99800000001

$ws.Range("H74").Value = 50450
$ws.Range("J74").Value = 50450
$ws.Range("L74").Value = 50450
$ws.Range("N74").Value = -52198

$ws.Range("H77").Value = 50450
$ws.Range("J77").Value = 50450
$ws.Range("L77").Value = 151350
$ws.Range("N77").Value = -160086

$ws.Range("H132").Value = 2441.0908
$ws.Range("I132").Value = 2298.2778
$ws.Range("K132").Value = 6894.8334
$ws.Range("M132").Value = -4364.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 91498.73
$ws.Range("I11").Value = 561.125
$ws.Range("J11").Value = 333999
$ws.Range("K11").Value = 1683.375
$ws.Range("L11").Value = 1001997
$ws.Range("M11").Value = -1543.375
$ws.Range("N11").Value = -1002277

$ws.Range("H14").Value = 502.2
$ws.Range("I14").Value = 502.2
$ws.Range("K14").Value = 1506.6
$ws.Range("M14").Value = -1333.6

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H29").Value = 333431.34
$ws.Range("I29").Value = 500074.5
$ws.Range("J29").Value = 145
$ws.Range("K29").Value = 1500223.5
$ws.Range("L29").Value = 435
$ws.Range("M29").Value = -1499946.5
$ws.Range("N29").Value = -989

$ws.Range("H32").Value = 334666
$ws.Range("I32").Value = 500499
$ws.Range("K32").Value = 1501497
$ws.Range("M32").Value = -1501214

$ws.Range("H34").Value = 1509.909
$ws.Range("I34").Value = 182
$ws.Range("J34").Value = 2616.5
$ws.Range("K34").Value = 546
$ws.Range("L34").Value = 7849.5
$ws.Range("M34").Value = -462
$ws.Range("N34").Value = -8017.5

$ws.Range("H46").Value = 91195.55
$ws.Range("I46").Value = 111410.11
$ws.Range("K46").Value = 334230.33
$ws.Range("M46").Value = -334139.33

$ws.Range("H62").Value = 150998.72
$ws.Range("I62").Value = 255998
$ws.Range("K62").Value = 767994
$ws.Range("M62").Value = -767308

$ws.Range("H65").Value = 150998.72
$ws.Range("I65").Value = 255998
$ws.Range("K65").Value = 2303982
$ws.Range("M65").Value = -2300550

$ws.Range("H107").Value = 614.2381
$ws.Range("I107").Value = 498.8889
$ws.Range("J107").Value = 700.75
$ws.Range("K107").Value = 1496.6667
$ws.Range("L107").Value = 2102.25
$ws.Range("M107").Value = 423.3333
$ws.Range("N107").Value = -5942.25

$ws.Range("H131").Value = 86038.664
$ws.Range("I131").Value = 334989
$ws.Range("K131").Value = 1004967
$ws.Range("M131").Value = -999927

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3356.5557
$ws.Range("I126").Value = 3262.2
$ws.Range("K126").Value = 9786.599999999999
$ws.Range("M126").Value = -7316.599999999999

$ws.Range("H132").Value = 5229.5
$ws.Range("I132").Value = 6369.1665
$ws.Range("J132").Value = 3520
$ws.Range("K132").Value = 19107.4995
$ws.Range("L132").Value = 10560
$ws.Range("M132").Value = -16577.4995
$ws.Range("N132").Value = -15620

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12999.818
$ws.Range("I122").Value = 16000
$ws.Range("K122").Value = 48000
$ws.Range("M122").Value = -45550

$ws.Range("H132").Value = 2943.7778
$ws.Range("I132").Value = 1857
$ws.Range("K132").Value = 5571
$ws.Range("M132").Value = -3041

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 25999.666
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 25999.666
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 25999.666
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -26981.666

$ws.Range("H122").Value = 2004.7407
$ws.Range("I122").Value = 1473.65
$ws.Range("K122").Value = 4420.950000000001
$ws.Range("M122").Value = -1970.950000000001
